$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles, fill, borders) of the last existing data row (row 32)
# down into the new rows (33-45) that we are about to populate.
$ws.Range("A32:E32").Copy()
$ws.Range("A33:E45").PasteSpecial(-4122)

# Match the row height used by all the other data rows.
for ($r = 33; $r -le 45; $r++) {
    $ws.Rows.Item($r).RowHeight = 20.1
}

# New i18n key/value rows: A = key, B = zh-CN, C = en-US
$newRows = @(
    @(33, "general.chooseFile", "选择文件", "Choose File"),
    @(34, "general.confirm", "确认", "Confirm"),
    @(35, "general.save", "保存", "Save"),
    @(36, "general.download", "下载", "Download"),
    @(37, "qrcode.basic.title", "基本", "Basic"),
    @(38, "qrcode.basic.placeholder", "链接或文本", "URL or Text"),
    @(39, "qrcode.basic.type", "类型", "Type"),
    @(40, "qrcode.basic.wifi", "WI-FI", "WI-FI"),
    @(41, "qrcode.basic.text", "文本", "Text"),
    @(42, "qrcode.advanced.title", "高级", "Advanced"),
    @(43, "qrcode.advanced.icon", "图标", "Icon"),
    @(44, "qrcode.advanced.light", "亮色", "Light Color"),
    @(45, "qrcode.advanced.dark", "暗色", "Dark Color"),
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

Write-Host "Dimension now: $($ws.UsedRange.Address())"
